# Fruta / hortaliza, semanal
# Insert a new weekly record at row 156, pushing the existing rows
# 156-180 down to 157-181 (the last existing row becomes row 181).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 156.
$ws.Rows("156:156").Insert()

# Populate the new row with the new weekly observation.
$ws.Cells.Item(156, 1).Value = 5
$ws.Cells.Item(156, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(156, 3).Value = "Maule"
$ws.Cells.Item(156, 4).Value = 44474
$ws.Cells.Item(156, 5).Value = 7
$ws.Cells.Item(156, 6).Value = 100114014
$ws.Cells.Item(156, 7).Value = "Betarraga"
$ws.Cells.Item(156, 8).Value = "Sin especificar"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 5000
$ws.Cells.Item(156, 11).Value = 700
$ws.Cells.Item(156, 12).Value = 700
$ws.Cells.Item(156, 13).Value = 700
$ws.Cells.Item(156, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(156, 15).Value = "Región del Maule"
$ws.Cells.Item(156, 16).Value = 140
$ws.Cells.Item(156, 17).Value = 5
$ws.Cells.Item(156, 18).Value = "Hortaliza"
